$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 45: "Relevant litteratur om systemtest" / "Any Role" / 2020-03-12 / 09:00-15:00
$ws.Range("A45").Value = "Relevant litteratur om systemtest"
$ws.Range("B45").Value = "Any Role"
$ws.Range("C45").Value = 43902
$ws.Range("D45").Value = 0.375
$ws.Range("E45").Value = 0.625
$ws.Range("G45").Formula = "=E45-D45"

# Row 46: "Opdatering af DD Dataordbog" / "business-Process Analyst" / 2020-03-13 / 16:40-17:15
$ws.Range("A46").Value = "Opdatering af DD Dataordbog"
$ws.Range("B46").Value = "business-Process Analyst"
$ws.Range("C46").Value = 43903
$ws.Range("D46").Value = 0.69444444444444453
$ws.Range("E46").Value = 0.71875
$ws.Range("G46").Formula = "=E46-D46"

# Update the on-disk selection / active cell to match the author's final cursor position
$ws.Range("B47").Select() | Out-Null
